$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns keep their text formatting
# (values like "1.00" / "0.999" would otherwise be auto-converted to numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '58.173.63'
$ws.Range("E2").Value = '  -1.28%  '
$ws.Range("D3").Value = '2.473.36'
$ws.Range("E3").Value = '  -1.95%  '
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").Value = '520.96'
$ws.Range("E5").Value = '  -2.97%  '
$ws.Range("D6").Value = '132.90'
$ws.Range("E6").Value = '  -3.73%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '0.559'
$ws.Range("E8").Value = '  -1.62%  '
$ws.Range("D9").Value = '0.0993'
$ws.Range("E9").Value = '  -2.20%  '
$ws.Range("E10").Value = '  -0.55%  '
$ws.Range("E11").Value = '  +0.35%  '
$ws.Range("D12").Value = '0.343'
$ws.Range("E12").Value = '  -1.77%  '
$ws.Range("D13").Value = '2.910.89'
$ws.Range("E13").Value = '  -2.08%  '
$ws.Range("D14").Value = '58.097.70'
$ws.Range("E14").Value = '  -1.47%  '
$ws.Range("D15").Value = '22.10'
$ws.Range("E15").Value = '  -4.52%  '
$ws.Range("D16").Value = '0.0000137'
$ws.Range("E16").Value = '  -2.18%  '
$ws.Range("D17").Value = '2.476.52'
$ws.Range("E17").Value = '  -1.88%  '
$ws.Range("D18").Value = '10.87'
$ws.Range("E18").Value = '  -2.30%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = '321.02'
$ws.Range("E19").Value = '  -1.47%  '
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").Value = '4.18'
$ws.Range("E20").Value = '  -2.70%  '
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("D22").Value = '5.76'
$ws.Range("E22").Value = '  -3.68%  '
$ws.Range("D23").Value = '64.49'
$ws.Range("E23").Value = '  -2.29%  '
$ws.Range("D24").Value = '0.409'
$ws.Range("E24").Value = '  -3.53%  '
$ws.Range("E25").Value = '  -0.13%  '
$ws.Range("E26").Value = '  -3.64%  '
$ws.Range("D27").Value = '7.42'
$ws.Range("E27").Value = '  -3.29%  '
$ws.Range("D28").Value = '0.0₃0753'
$ws.Range("E28").Value = '  -2.85%  '
$ws.Range("D29").Value = '6.40'
$ws.Range("E29").Value = '  -4.89%  '
$ws.Range("D30").Value = '1.71'
$ws.Range("E30").Value = '  -4.73%  '
$ws.Range("D31").Value = '165.39'
$ws.Range("E31").Value = '  +1.88%  '
$ws.Range("E32").Value = '  -4.02%  '
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("E34").Value = '  +0.11%  '
$ws.Range("D35").Value = '18.17'
$ws.Range("E35").Value = '  -1.83%  '
$ws.Range("E36").Value = '  -10.12%  '
$ws.Range("D37").Value = '4.00'
$ws.Range("E37").Value = '  -3.23%  '
$ws.Range("E38").Value = '  -3.67%  '
$ws.Range("D39").Value = '0.799'
$ws.Range("E39").Value = '  -2.40%  '
$ws.Range("D40").Value = '276.94'
$ws.Range("E40").Value = '  -3.43%  '
$ws.Range("D41").Value = '3.48'
$ws.Range("E41").Value = '  -4.65%  '
$ws.Range("D42").Value = '5.05'
$ws.Range("E42").Value = '  -3.33%  '
$ws.Range("D43").Value = '0.596'
$ws.Range("E43").Value = '  -2.54%  '
$ws.Range("D44").Value = '126.14'
$ws.Range("E44").Value = '  -4.85%  '
$ws.Range("D45").Value = '0.0909'
$ws.Range("D46").Value = '0.0494'
$ws.Range("E46").Value = '  -3.40%  '
$ws.Range("D47").Value = '0.0215'
$ws.Range("E47").Value = '  -3.32%  '
$ws.Range("D48").Value = '17.16'
$ws.Range("D49").Value = '1.740.34'
$ws.Range("E49").Value = '  -1.26%  '
$ws.Range("D50").Value = '0.973'
$ws.Range("E50").Value = '  -1.53%  '
$ws.Range("D51").Value = '4.68'
$ws.Range("E51").Value = '  -1.83%  '

# Restore default styling on the touched columns (no visual/style change intended).
$ws.Range("D2:E51").Style = "Normal"

